# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same as the existing header row (bold,
# bordered, centered) by copying the format from the last existing header
# cell (AC1) before overwriting the values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-52): same W/L/T totals
# repeated down the whole column, as plain numbers (no special style).
$ws.Range("AD2:AD52").Value = 64
$ws.Range("AE2:AE52").Value = 98
$ws.Range("AF2:AF52").Value = 0

Write-Output "season record columns added"
